$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings in column D stay as text (matches source data: prices
# are stored as literal strings, e.g. "45.227.25" / "0.0800", not numbers).
$priceCells = @("D2","D3","D4","D5","D6","D9","D10","D11","D13","D14","D15","D16","D17","D18","D19","D20","D21","D22","D23","D24","D25","D27","D28","D30","D31","D32","D33","D36","D37","D38","D39","D40","D41","D43","D45","D47","D48","D49","D50")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "45.227.25"
$ws.Range("E2").Value = "  +1.63%  "
$ws.Range("D3").Value = "2.421.85"
$ws.Range("E3").Value = "  -0.36%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "318.66"
$ws.Range("E5").Value = "  +3.14%  "
$ws.Range("D6").Value = "103.18"
$ws.Range("E6").Value = "  +1.57%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("D9").Value = "0.532"
$ws.Range("E9").Value = "  +4.82%  "
$ws.Range("D10").Value = "35.55"
$ws.Range("E10").Value = "  +0.97%  "
$ws.Range("D11").Value = "0.0800"
$ws.Range("E11").Value = "  -0.43%  "
$ws.Range("E12").Value = "  -1.43%  "
$ws.Range("D13").Value = "18.11"
$ws.Range("E13").Value = "  -3.64%  "
$ws.Range("D14").Value = "6.98"
$ws.Range("E14").Value = "  +0.05%  "
$ws.Range("D15").Value = "2.799.58"
$ws.Range("E15").Value = "  -0.41%  "
$ws.Range("D16").Value = "2.400.01"
$ws.Range("E16").Value = "  -2.10%  "
$ws.Range("D17").Value = "0.834"
$ws.Range("E17").Value = "  -0.49%  "
$ws.Range("D18").Value = "45.157.33"
$ws.Range("E18").Value = "  +1.58%  "
$ws.Range("D19").Value = "12.19"
$ws.Range("E19").Value = "  -2.78%  "
$ws.Range("D20").Value = "6.33"
$ws.Range("E20").Value = "  -1.12%  "
$ws.Range("D21").Value = "0.0₃0925"
$ws.Range("E21").Value = "  +1.77%  "
$ws.Range("D22").Value = "70.42"
$ws.Range("E22").Value = "  +2.28%  "
$ws.Range("B23").Value = "ImmutableX"
$ws.Range("C23").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D23").Value = "2.36"
$ws.Range("E23").Value = "  +0.81%  "
$ws.Range("B24").Value = "BitcoinCash"
$ws.Range("C24").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D24").Value = "245.34"
$ws.Range("E24").Value = "  +1.34%  "
$ws.Range("D25").Value = "2.48"
$ws.Range("E25").Value = "  -0.35%  "
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("D27").Value = "25.60"
$ws.Range("E27").Value = "  +1.17%  "
$ws.Range("D28").Value = "2.28"
$ws.Range("E28").Value = "  +4.54%  "
$ws.Range("E29").Value = "  -0.35%  "
$ws.Range("B30").Value = "InjectiveProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D30").Value = "33.08"
$ws.Range("E30").Value = "  -1.03%  "
$ws.Range("B31").Value = "OKB"
$ws.Range("C31").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D31").Value = "48.95"
$ws.Range("E31").Value = "  +0.36%  "
$ws.Range("D32").Value = "20.49"
$ws.Range("E32").Value = "  +5.79%  "
$ws.Range("D33").Value = "0.127"
$ws.Range("E33").Value = "  +6.76%  "
$ws.Range("E34").Value = "  -0.03%  "
$ws.Range("E35").Value = "  +0.30%  "
$ws.Range("D36").Value = "0.0757"
$ws.Range("E36").Value = "  -1.00%  "
$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D37").Value = "4.50"
$ws.Range("E37").Value = "  -0.20%  "
$ws.Range("B38").Value = "ARBITRUM"
$ws.Range("C38").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D38").Value = "1.87"
$ws.Range("E38").Value = "  -1.88%  "
$ws.Range("B39").Value = "LidoDAOToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D39").Value = "2.90"
$ws.Range("E39").Value = "  -0.38%  "
$ws.Range("B40").Value = "Monero"
$ws.Range("C40").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D40").Value = "126.37"
$ws.Range("E40").Value = "  -4.55%  "
$ws.Range("D41").Value = "2.28"
$ws.Range("E41").Value = "  -2.83%  "
$ws.Range("E42").Value = "  +0.54%  "
$ws.Range("D43").Value = "20.72"
$ws.Range("E43").Value = "  -5.53%  "
$ws.Range("E44").Value = "  -0.21%  "
$ws.Range("D45").Value = "1.941.30"
$ws.Range("E45").Value = "  -1.05%  "
$ws.Range("E46").Value = "  -2.73%  "
$ws.Range("D47").Value = "2.94"
$ws.Range("E47").Value = "  +0.69%  "
$ws.Range("B48").Value = "FraxShare"
$ws.Range("C48").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D48").Value = "9.14"
$ws.Range("E48").Value = "  -3.28%  "
$ws.Range("B49").Value = "Stacks"
$ws.Range("C49").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D49").Value = "1.79"
$ws.Range("E49").Value = "  +5.95%  "
$ws.Range("D50").Value = "78.33"
$ws.Range("E50").Value = "  +6.01%  "
$ws.Range("E51").Value = "  +3.79%  "
